$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.020402789115906
$ws.Range("B1").Value = 2.20635461807251
$ws.Range("C1").Value = 7.109145164489746
$ws.Range("D1").Value = 2.365393400192261
$ws.Range("E1").Value = 1.331662893295288
